{"js": "// Apply the \"2024-08-06 Tuesday\" worksheet update: the header date and all\n// 25 multiplication-answer cells are replaced with new values (same count,\n// same positions - only the text content changes).\nconst replacements = [\n  [\"2024-08-05 Monday\", \"2024-08-06 Tuesday\"],\n  [\"54\u00d711=594\", \"92\u00d744=4048\"],\n  [\"37\u00d737=1369\", \"20\u00d732=640\"],\n  [\"15\u00d794=1410\", \"65\u00d722=1430\"],\n  [\"16\u00d731=496\", \"50\u00d779=3950\"],\n  [\"23\u00d765=1495\", \"57\u00d743=2451\"],\n  [\"94\u00d791=8554\", \"26\u00d754=1404\"],\n  [\"58\u00d780=4640\", \"64\u00d764=4096\"],\n  [\"89\u00d727=2403\", \"61\u00d771=4331\"],\n  [\"32\u00d730=960\", \"97\u00d772=6984\"],\n  [\"82\u00d714=1148\", \"38\u00d797=3686\"],\n  [\"43\u00d754=2322\", \"53\u00d751=2703\"],\n  [\"70\u00d743=3010\", \"76\u00d745=3420\"],\n  [\"50\u00d780=4000\", \"67\u00d771=4757\"],\n  [\"97\u00d736=3492\", \"92\u00d727=2484\"],\n  [\"61\u00d787=5307\", \"45\u00d773=3285\"],\n  [\"37\u00d740=1480\", \"78\u00d773=5694\"],\n  [\"54\u00d763=3402\", \"20\u00d717=340\"],\n  [\"34\u00d725=850\", \"63\u00d798=6174\"],\n  [\"35\u00d744=1540\", \"32\u00d779=2528\"],\n  [\"78\u00d769=5382\", \"29\u00d793=2697\"],\n  [\"44\u00d760=2640\", \"79\u00d771=5609\"],\n  [\"77\u00d743=3311\", \"63\u00d761=3843\"],\n  [\"98\u00d766=6468\", \"48\u00d773=3504\"],\n  [\"36\u00d796=3456\", \"43\u00d785=3655\"],\n  [\"85\u00d774=6290\", \"75\u00d796=7200\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"2024-08-06 Tuesday\" worksheet update: the header date and all\n# 25 multiplication-answer cells are replaced with new values (same count,\n# same positions - only the text content changes).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-05 Monday\", \"2024-08-06 Tuesday\"),\n    @(\"54\u00d711=594\", \"92\u00d744=4048\"),\n    @(\"37\u00d737=1369\", \"20\u00d732=640\"),\n    @(\"15\u00d794=1410\", \"65\u00d722=1430\"),\n    @(\"16\u00d731=496\", \"50\u00d779=3950\"),\n    @(\"23\u00d765=1495\", \"57\u00d743=2451\"),\n    @(\"94\u00d791=8554\", \"26\u00d754=1404\"),\n    @(\"58\u00d780=4640\", \"64\u00d764=4096\"),\n    @(\"89\u00d727=2403\", \"61\u00d771=4331\"),\n    @(\"32\u00d730=960\", \"97\u00d772=6984\"),\n    @(\"82\u00d714=1148\", \"38\u00d797=3686\"),\n    @(\"43\u00d754=2322\", \"53\u00d751=2703\"),\n    @(\"70\u00d743=3010\", \"76\u00d745=3420\"),\n    @(\"50\u00d780=4000\", \"67\u00d771=4757\"),\n    @(\"97\u00d736=3492\", \"92\u00d727=2484\"),\n    @(\"61\u00d787=5307\", \"45\u00d773=3285\"),\n    @(\"37\u00d740=1480\", \"78\u00d773=5694\"),\n    @(\"54\u00d763=3402\", \"20\u00d717=340\"),\n    @(\"34\u00d725=850\", \"63\u00d798=6174\"),\n    @(\"35\u00d744=1540\", \"32\u00d779=2528\"),\n    @(\"78\u00d769=5382\", \"29\u00d793=2697\"),\n    @(\"44\u00d760=2640\", \"79\u00d771=5609\"),\n    @(\"77\u00d743=3311\", \"63\u00d761=3843\"),\n    @(\"98\u00d766=6468\", \"48\u00d773=3504\"),\n    @(\"36\u00d796=3456\", \"43\u00d785=3655\"),\n    @(\"85\u00d774=6290\", \"75\u00d796=7200\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
